$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "GATTOO"
$ws.Range("B6").Value = "GATTOO"
$ws.Range("C6").Value = $true

$ws.Range("C6").Select()
